# Update generated statistics (column F) across the four sheets, per the
# "Update gh-pages to output generated at 456a3b4" data refresh.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1)
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2,6).Value  = 12804
$ws.Cells.Item(3,6).Value  = 7163
$ws.Cells.Item(6,6).Value  = 450
$ws.Cells.Item(10,6).Value = 1005
$ws.Cells.Item(11,6).Value = 145
$ws.Cells.Item(13,6).Value = 1013
$ws.Cells.Item(18,6).Value = 246
$ws.Cells.Item(22,6).Value = 311
$ws.Cells.Item(24,6).Value = 162
$ws.Cells.Item(26,6).Value = 5234
$ws.Cells.Item(29,6).Value = 311
$ws.Cells.Item(30,6).Value = 1361
$ws.Cells.Item(32,6).Value = 39
$ws.Cells.Item(36,6).Value = 594

# Sheet "演出" (index 2)
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(4,6).Value  = 3743
$ws.Cells.Item(5,6).Value  = 3743
$ws.Cells.Item(19,6).Value = 20

# Sheet "本地生活" (index 3)
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2,6).Value = 9279
$ws.Cells.Item(4,6).Value = 2014

# Sheet "全部类型" (index 4)
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2,6).Value  = 9279
$ws.Cells.Item(4,6).Value  = 2014
$ws.Cells.Item(5,6).Value  = 12804
$ws.Cells.Item(6,6).Value  = 7163
$ws.Cells.Item(8,6).Value  = 3743
$ws.Cells.Item(9,6).Value  = 450
$ws.Cells.Item(10,6).Value = 1005
$ws.Cells.Item(11,6).Value = 145
$ws.Cells.Item(13,6).Value = 1013
$ws.Cells.Item(18,6).Value = 246
$ws.Cells.Item(22,6).Value = 311
$ws.Cells.Item(27,6).Value = 162
$ws.Cells.Item(29,6).Value = 5234
$ws.Cells.Item(34,6).Value = 311
$ws.Cells.Item(36,6).Value = 1361
$ws.Cells.Item(40,6).Value = 594
$ws.Cells.Item(48,6).Value = 20
